$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture an unstyled "plain" style to reuse so forcing text format on column D
# cells does not leave a stray custom style applied to the cell.
$plainStyle = $ws.Range("C2").Style

# --- Column D (Price) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.436.19"
$ws.Range("D2").Style = $plainStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.081.72"
$ws.Range("D3").Style = $plainStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.55"
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = $plainStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.24"
$ws.Range("D7").Style = $plainStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("D9").Style = $plainStyle
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0833"
$ws.Range("D10").Style = $plainStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = $plainStyle
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.389.15"
$ws.Range("D12").Style = $plainStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.87"
$ws.Range("D13").Style = $plainStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.39"
$ws.Range("D14").Style = $plainStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("D15").Style = $plainStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.44"
$ws.Range("D16").Style = $plainStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.080.30"
$ws.Range("D17").Style = $plainStyle
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.364.70"
$ws.Range("D18").Style = $plainStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.35"
$ws.Range("D19").Style = $plainStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.00"
$ws.Range("D20").Style = $plainStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.88"
$ws.Range("D22").Style = $plainStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = $plainStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.10"
$ws.Range("D26").Style = $plainStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("D27").Style = $plainStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("D28").Style = $plainStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("D29").Style = $plainStyle
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = $plainStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("D31").Style = $plainStyle
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("D34").Style = $plainStyle
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0605"
$ws.Range("D35").Style = $plainStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.35"
$ws.Range("D37").Style = $plainStyle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("D38").Style = $plainStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = $plainStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.20"
$ws.Range("D40").Style = $plainStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.542.57"
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.31"
$ws.Range("D42").Style = $plainStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("D43").Style = $plainStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0921"
$ws.Range("D44").Style = $plainStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.62"
$ws.Range("D46").Style = $plainStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.276.59"
$ws.Range("D51").Style = $plainStyle

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("E14").Value = "  +6.70%  "
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +6.98%  "
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("E30").Value = "  +8.68%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("E33").Value = "  +7.26%  "
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +7.65%  "
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  +2.41%  "

